$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.427.10'
$ws.Range('E2').Value = '  +1.76%  '
$ws.Range('D3').Value = '2.156.11'
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.95'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.03%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +2.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0855'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.94'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.50%  '
$ws.Range('D13').Value = '2.477.49'
$ws.Range('E13').Value = '  +2.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.811'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').Value = '2.176.04'
$ws.Range('E17').Value = '  +3.81%  '
$ws.Range('D18').Value = '39.352.98'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.65'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.37%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '172.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.141'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.86'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  +7.02%  '
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('E34').Value = '  +8.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.74'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '103.55'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('D42').Value = '1.535.97'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.47%  '
$ws.Range('E45').Value = '  +2.34%  '
$ws.Range('E46').Value = '  +4.53%  '
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('E48').Value = '  +4.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.67'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('D50').Value = '2.360.93'
$ws.Range('E50').Value = '  +2.98%  '
$ws.Range('E51').Value = '  +0.00%  '
